$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 60
$ws.Cells.Item(60, 1).Value = 112103528
$ws.Cells.Item(60, 2).Value = 8377
$ws.Cells.Item(60, 4).Value = "LC"
$ws.Cells.Item(60, 5).Value = 106545
$ws.Cells.Item(60, 6).Value = "Mindre märgborre"
$ws.Cells.Item(60, 7).Value = "Tomicus minor"
$ws.Cells.Item(60, 8).Value = "(Hartig, 1834)"
$ws.Cells.Item(60, 11).ClearContents()
$ws.Cells.Item(60, 13).Value = "äldre gnagspår"
$ws.Cells.Item(60, 17).Value = 572399.2401777974
$ws.Cells.Item(60, 18).Value = 6634971.688825586

# Row 61
$ws.Cells.Item(61, 1).Value = 112103561
$ws.Cells.Item(61, 2).Value = 5113
$ws.Cells.Item(61, 4).Value = "LC"
$ws.Cells.Item(61, 5).Value = 100526
$ws.Cells.Item(61, 6).Value = "Bronshjon"
$ws.Cells.Item(61, 7).Value = "Callidium coriaceum"
$ws.Cells.Item(61, 8).Value = "Paykull, 1800"
$ws.Cells.Item(61, 11).ClearContents()
$ws.Cells.Item(61, 13).Value = "färska gnagspår"
$ws.Cells.Item(61, 17).Value = 572285.1631843462
$ws.Cells.Item(61, 18).Value = 6634908.643639773

# Row 62
$ws.Cells.Item(62, 1).Value = 112103532
$ws.Cells.Item(62, 2).Value = 96348
$ws.Cells.Item(62, 4).Value = "VU"
$ws.Cells.Item(62, 5).Value = 220787
$ws.Cells.Item(62, 6).Value = "Knärot"
$ws.Cells.Item(62, 7).Value = "Goodyera repens"
$ws.Cells.Item(62, 8).Value = "(L.) R. Br."
$ws.Cells.Item(62, 11).Value = "fullt utvecklade blad"
$ws.Cells.Item(62, 13).ClearContents()
$ws.Cells.Item(62, 17).Value = 572381.684763086
$ws.Cells.Item(62, 18).Value = 6635329.725023665

# Row 63
$ws.Cells.Item(63, 1).Value = 112103548
$ws.Cells.Item(63, 2).Value = 96348
$ws.Cells.Item(63, 4).Value = "VU"
$ws.Cells.Item(63, 5).Value = 220787
$ws.Cells.Item(63, 6).Value = "Knärot"
$ws.Cells.Item(63, 7).Value = "Goodyera repens"
$ws.Cells.Item(63, 8).Value = "(L.) R. Br."
$ws.Cells.Item(63, 11).Value = "fullt utvecklade blad"
$ws.Cells.Item(63, 13).ClearContents()
$ws.Cells.Item(63, 17).Value = 572360.7400891574
$ws.Cells.Item(63, 18).Value = 6635164.45678684

# Row 64
$ws.Cells.Item(64, 1).Value = 112103517
$ws.Cells.Item(64, 2).Value = 89405
$ws.Cells.Item(64, 4).Value = "NT"
$ws.Cells.Item(64, 5).Value = 1202
$ws.Cells.Item(64, 6).Value = "Ullticka"
$ws.Cells.Item(64, 7).Value = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(64, 8).Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(64, 11).ClearContents()
$ws.Cells.Item(64, 13).ClearContents()
$ws.Cells.Item(64, 17).Value = 572367.7229917983
$ws.Cells.Item(64, 18).Value = 6634935.385356643

# Row 65
$ws.Cells.Item(65, 1).Value = 112103511
$ws.Cells.Item(65, 2).Value = 94134
$ws.Cells.Item(65, 4).Value = "NT"
$ws.Cells.Item(65, 5).Value = 53
$ws.Cells.Item(65, 6).Value = "Vedtrappmossa"
$ws.Cells.Item(65, 7).Value = "Crossocalyx hellerianus"
$ws.Cells.Item(65, 8).Value = "(Nees ex Lindenb.) Meyl."
$ws.Cells.Item(65, 11).ClearContents()
$ws.Cells.Item(65, 13).ClearContents()
$ws.Cells.Item(65, 17).Value = 572312.6870405492
$ws.Cells.Item(65, 18).Value = 6634839.810999912

# Row 66
$ws.Cells.Item(66, 1).Value = 112103529
$ws.Cells.Item(66, 2).Value = 8377
$ws.Cells.Item(66, 4).Value = "LC"
$ws.Cells.Item(66, 5).Value = 106545
$ws.Cells.Item(66, 6).Value = "Mindre märgborre"
$ws.Cells.Item(66, 7).Value = "Tomicus minor"
$ws.Cells.Item(66, 8).Value = "(Hartig, 1834)"
$ws.Cells.Item(66, 11).ClearContents()
$ws.Cells.Item(66, 13).Value = "äldre gnagspår"
$ws.Cells.Item(66, 17).Value = 572285.7687057387
$ws.Cells.Item(66, 18).Value = 6634929.264707729

# Row 73
$ws.Cells.Item(73, 1).Value = 112103542
$ws.Cells.Item(73, 2).Value = 96348
$ws.Cells.Item(73, 4).Value = "VU"
$ws.Cells.Item(73, 5).Value = 220787
$ws.Cells.Item(73, 6).Value = "Knärot"
$ws.Cells.Item(73, 7).Value = "Goodyera repens"
$ws.Cells.Item(73, 8).Value = "(L.) R. Br."
$ws.Cells.Item(73, 11).Value = "fullt utvecklade blad"
$ws.Cells.Item(73, 13).ClearContents()
$ws.Cells.Item(73, 17).Value = 572414.0757496187
$ws.Cells.Item(73, 18).Value = 6635062.958343645

# Row 74
$ws.Cells.Item(74, 1).Value = 112103545
$ws.Cells.Item(74, 2).Value = 96348
$ws.Cells.Item(74, 4).Value = "VU"
$ws.Cells.Item(74, 5).Value = 220787
$ws.Cells.Item(74, 6).Value = "Knärot"
$ws.Cells.Item(74, 7).Value = "Goodyera repens"
$ws.Cells.Item(74, 8).Value = "(L.) R. Br."
$ws.Cells.Item(74, 11).Value = "fullt utvecklade blad"
$ws.Cells.Item(74, 13).ClearContents()
$ws.Cells.Item(74, 17).Value = 572405.2016443094
$ws.Cells.Item(74, 18).Value = 6634975.826361955

# Row 75
$ws.Cells.Item(75, 1).Value = 112103520
$ws.Cells.Item(75, 2).Value = 89425
$ws.Cells.Item(75, 4).Value = "NT"
$ws.Cells.Item(75, 5).Value = 5442
$ws.Cells.Item(75, 6).Value = "Tallticka"
$ws.Cells.Item(75, 7).Value = "Porodaedalea pini"
$ws.Cells.Item(75, 8).Value = "(Brot.) Murrill"
$ws.Cells.Item(75, 11).ClearContents()
$ws.Cells.Item(75, 13).ClearContents()
$ws.Cells.Item(75, 17).Value = 572409.3458258022
$ws.Cells.Item(75, 18).Value = 6634969.875286552

# Row 76
$ws.Cells.Item(76, 1).Value = 112103539
$ws.Cells.Item(76, 2).Value = 96348
$ws.Cells.Item(76, 4).Value = "VU"
$ws.Cells.Item(76, 5).Value = 220787
$ws.Cells.Item(76, 6).Value = "Knärot"
$ws.Cells.Item(76, 7).Value = "Goodyera repens"
$ws.Cells.Item(76, 8).Value = "(L.) R. Br."
$ws.Cells.Item(76, 11).Value = "fullt utvecklade blad"
$ws.Cells.Item(76, 13).ClearContents()
$ws.Cells.Item(76, 17).Value = 572444.8101381793
$ws.Cells.Item(76, 18).Value = 6635165.091275458
